# unify the conception of DataNode, DataTable, Entity.
# Rename the two worksheets and make the second one ("DataTable") the
# active/selected sheet, matching the authored workbook.

$wb = $excel.ActiveWorkbook

$wsNode = $wb.Worksheets.Item(1)
$wsTable = $wb.Worksheets.Item(2)

$wsNode.Name = "DataNode"
$wsTable.Name = "DataTable"

# Make "DataTable" the active tab (was "Property1"/DataNode before).
$wsTable.Activate()
